$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the broken-path folder name ("Python Learning" -> "Python-Learning")
# that was causing PyInstaller to not find the binary, and populate the
# output_file_path / output_file_name columns that the script now fills in.

$ws.Range("A2").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\1_input\SBM_MV_1080×1080_220117_iP13pro_gp_nts_YDN.mp4"
$ws.Range("D2").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\2_output\SBM_MV_1080×1080_220117_iP13pro_gp_nts_YDN.jpg"
$ws.Range("E2").Value = "SBM_MV_1080×1080_220117_iP13pro_gp_nts_YDN.jpg"

$ws.Range("A3").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\1_input\SBM_MV_1080×1080_220117_iP13_bl_olaf_YDN.mp4"
$ws.Range("D3").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\2_output\SBM_MV_1080×1080_220117_iP13_bl_olaf_YDN.jpg"
$ws.Range("E3").Value = "SBM_MV_1080×1080_220117_iP13_bl_olaf_YDN.jpg"

$ws.Range("A4").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\1_input\SBM_MV_1280×720_220117_iP13pro_sibl_nts_YDN.mp4"
$ws.Range("D4").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\2_output\SBM_MV_1280×720_220117_iP13pro_sibl_nts_YDN.jpg"
$ws.Range("E4").Value = "SBM_MV_1280×720_220117_iP13pro_sibl_nts_YDN.jpg"

$ws.Range("A5").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\1_input\SBM_MV_1280×720_220117_iP13_mdn_olaf_YDN.mp4"
$ws.Range("D5").Value = "C:\Users\DT0083\Desktop\Python-Learning\ezVideoThumbnails\2_output\SBM_MV_1280×720_220117_iP13_mdn_olaf_YDN.jpg"
$ws.Range("E5").Value = "SBM_MV_1280×720_220117_iP13_mdn_olaf_YDN.jpg"
